# Correction s_vs_i valeur temporelle de s+ et s-
#
# Inserts two new computed columns (s+ [m] and s-[mm], in F and G) ahead of
# the existing "s_vs_i" block in the second table (row 12 headers, rows
# 13-19 data), shifting the old F:O block to H:Q, and adds the new
# "distance parcourue par le moteur" constant used by the new formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New header cells next to the first table (row 1) + the new distance
#    constant referenced by the new formulas as $K$2.
# ---------------------------------------------------------------------
$ws.Range("J1").Value = " "
$ws.Range("K1").Value = "distance parcourue par le moteur :"
$ws.Range("K2").Value = 25.5524

# ---------------------------------------------------------------------
# 2) Shift the second table's data block (rows 13-19) two columns to the
#    right: old F..O -> new H..Q. Work from the rightmost column back to
#    the leftmost so that source cells are read before they're
#    overwritten by another column's move.
# ---------------------------------------------------------------------
for ($r = 13; $r -le 19; $r++) {

    # O -> Q (value)
    $ws.Cells.Item($r, 17).Value = $ws.Cells.Item($r, 15).Value2
    # N -> P (value)
    $ws.Cells.Item($r, 16).Value = $ws.Cells.Item($r, 14).Value2
    # M -> O (formula: (K*J)/I -> (M*L)/K)
    $ws.Cells.Item($r, 15).Formula = "=(M$r*L$r)/K$r"
    # L -> N (formula: ((0.5*25.5524)/(G-F)) -> ((0.5*25.5524)/(I-H)))
    $ws.Cells.Item($r, 14).Formula = "=((0.5*25.5524)/(I$r-H$r))"
    # K -> M (formula: (L*I)/H -> (N*K)/J)
    $ws.Cells.Item($r, 13).Formula = "=(N$r*K$r)/J$r"
    # J -> L (formula: 5.5*0.1*I -> 5.5*0.1*K)
    $ws.Cells.Item($r, 12).Formula = "=5.5 *0.1* K$r"
    # I -> K (formula text is a literal expression, unchanged)
    $iFormula = $ws.Cells.Item($r, 9).Formula
    $ws.Cells.Item($r, 11).Formula = $iFormula
    # H -> J (formula: ((D-E)/2*25.5524)/(G-F) -> ((D-E)/2*25.5524)/(I-H))
    $ws.Cells.Item($r, 10).Formula = "=(((D$r-E$r)/2)*25.5524)/(I$r-H$r)"
    # G -> I (value)
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($r, 7).Value2
    # F -> H (value)
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($r, 6).Value2

    # New F, G formulas (computed from the shifted columns + $K$2)
    $ws.Cells.Item($r, 6).Formula = "=(D$r * `$K`$2)/(I$r-H$r)"
    $ws.Cells.Item($r, 7).Formula = "=(Q$r * `$K`$2)/(I$r-H$r)"
}

# ---------------------------------------------------------------------
# 3) Rewrite the second table's header row (row 12) for the new column
#    layout. Labels that simply move to a new column are carried over
#    (rightmost-first, to avoid clobbering a value before it is read);
#    labels that are brand new are typed in directly. The two old labels
#    "s [m]" (H12) and "Incertitude s [m]" (L12) are retired entirely.
# ---------------------------------------------------------------------
$ws.Range("Q12").Value = $ws.Range("O12").Value2   # |s-| [s]
$ws.Range("P12").Value = $ws.Range("N12").Value2   # Incertitude s [s]
$ws.Range("O12").Value = $ws.Range("M12").Value2   # Incertitude B moyen
$ws.Range("M12").Value = $ws.Range("K12").Value2   # Incertitude dB/dz
$ws.Range("L12").Value = $ws.Range("J12").Value2   # B moyen [T]
$ws.Range("K12").Value = $ws.Range("I12").Value2   # dB / dz [T/cm]
$ws.Range("I12").Value = $ws.Range("G12").Value2   # t0'[s]
$ws.Range("H12").Value = $ws.Range("F12").Value2   # t0 [s]

$ws.Range("F12").Value = "s+ [m]"
$ws.Range("G12").Value = "s-[mm]"
$ws.Range("J12").Value = "s [mm]"
$ws.Range("N12").Value = "Incertitude s [mm]"

# ---------------------------------------------------------------------
# 4) Cosmetic touch-ups matching the rest of the diff: widen column K,
#    size the two new columns, and move the selection to G13.
# ---------------------------------------------------------------------
$ws.Columns.Item(11).ColumnWidth = 30.15
$ws.Columns.Item(15).ColumnWidth = 16.45
$ws.Columns.Item(16).ColumnWidth = 16.6

$ws.Range("G13").Select()
